$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 555100
$ws.Range("E8").Value = 384900
$ws.Range("F8").Value = 258700
$ws.Range("G8").Value = 158700
$ws.Range("H8").Value = 90300
$ws.Range("I8").Value = 47500
$ws.Range("D9").Value = 464300
$ws.Range("E9").Value = 325600
$ws.Range("F9").Value = 226500
$ws.Range("G9").Value = 144300
$ws.Range("H9").Value = 80000
$ws.Range("I9").Value = 41200
$ws.Range("D10").Value = 90800
$ws.Range("E10").Value = 59400
$ws.Range("F10").Value = 32200
$ws.Range("G10").Value = 14400
$ws.Range("H10").Value = 10300
$ws.Range("I10").Value = 6300
$ws.Range("D12").Value = 9200
$ws.Range("E12").Value = 8100
$ws.Range("F12").Value = 6100
$ws.Range("G12").Value = 3700
$ws.Range("H12").Value = 1300
$ws.Range("D17").Value = 541000
$ws.Range("E17").Value = 389300
$ws.Range("F17").Value = 290100
$ws.Range("G17").Value = 176300
$ws.Range("H17").Value = 94600
$ws.Range("I17").Value = 50000
$ws.Range("D18").Value = 14100
$ws.Range("F18").Value = -31400
$ws.Range("G18").Value = -17500
$ws.Range("H18").Value = -4300
$ws.Range("E20").Value = -2300
$ws.Range("D21").Value = 17100
$ws.Range("E21").Value = -4600
$ws.Range("F21").Value = -31300
$ws.Range("G21").Value = -16600
$ws.Range("D23").Value = 15100
$ws.Range("E23").Value = -6600
$ws.Range("F23").Value = -32900
$ws.Range("G23").Value = -17700
$ws.Range("H23").Value = -4300
$ws.Range("D24").Value = -4700
$ws.Range("D26").Value = 19800
$ws.Range("E26").Value = -6600
$ws.Range("F26").Value = -32900
$ws.Range("G26").Value = -17700
$ws.Range("H26").Value = -4300
$ws.Range("D27").Value = -10300
$ws.Range("E27").Value = -95000
$ws.Range("F27").Value = -64700
$ws.Range("G27").Value = -34500
$ws.Range("H27").Value = -6400
$ws.Range("I27").Value = -3600
$ws.Range("E32").Value = 2300
$ws.Range("D33").Value = -10300
$ws.Range("E33").Value = -95000
$ws.Range("F33").Value = -64700
$ws.Range("G33").Value = -34500
$ws.Range("H33").Value = -6400
$ws.Range("I33").Value = -3600
$ws.Range("D35").Value = -10300
$ws.Range("E35").Value = -95000
$ws.Range("F35").Value = -64700
$ws.Range("G35").Value = -34500
$ws.Range("H35").Value = -6400
$ws.Range("I35").Value = -3600
$ws.Range("D41").Value = 67300
$ws.Range("E41").Value = 8200
$ws.Range("F41").Value = 42200
$ws.Range("G41").Value = 10700
$ws.Range("H41").Value = 7500
$ws.Range("D43").Value = 8100
$ws.Range("G43").Value = 2000
$ws.Range("I43").Value = 400
$ws.Range("D44").Value = 176600
$ws.Range("E44").Value = 111600
$ws.Range("F44").Value = 68900
$ws.Range("G44").Value = 80100
$ws.Range("H44").Value = 52600
$ws.Range("I44").Value = 22700
$ws.Range("D45").Value = 62900
$ws.Range("E45").Value = 26700
$ws.Range("F45").Value = 4600
$ws.Range("G45").Value = 19800
$ws.Range("D46").Value = 314800
$ws.Range("E46").Value = 149700
$ws.Range("F46").Value = 116900
$ws.Range("G46").Value = 112600
$ws.Range("H46").Value = 63100
$ws.Range("I46").Value = 26800
$ws.Range("D48").Value = 6100
$ws.Range("E48").Value = 5200
$ws.Range("F48").Value = 5300
$ws.Range("G48").Value = 5000
$ws.Range("D52").Value = 26100
$ws.Range("F52").Value = 23700
$ws.Range("H52").Value = 1200
$ws.Range("I52").Value = 800
$ws.Range("D54").Value = 346900
$ws.Range("E54").Value = 155200
$ws.Range("F54").Value = 145900
$ws.Range("G54").Value = 118400
$ws.Range("H54").Value = 66300
$ws.Range("I54").Value = 29700
$ws.Range("D57").Value = 47300
$ws.Range("E57").Value = 40800
$ws.Range("F57").Value = 42900
$ws.Range("G57").Value = 53100
$ws.Range("H57").Value = 42300
$ws.Range("I57").Value = 18400
$ws.Range("D58").Value = 26300
$ws.Range("E58").Value = 29700
$ws.Range("F58").Value = 26100
$ws.Range("G58").Value = 13400
$ws.Range("D59").Value = 63400
$ws.Range("E59").Value = 39300
$ws.Range("F59").Value = 29700
$ws.Range("G59").Value = 24400
$ws.Range("H59").Value = 9800
$ws.Range("D60").Value = 137000
$ws.Range("E60").Value = 109700
$ws.Range("F60").Value = 98800
$ws.Range("G60").Value = 90900
$ws.Range("H60").Value = 52900
$ws.Range("I60").Value = 22000
$ws.Range("D61").Value = 18500
$ws.Range("D66").Value = 156500
$ws.Range("E66").Value = 110800
$ws.Range("F66").Value = 98800
$ws.Range("G66").Value = 90900
$ws.Range("H66").Value = 52900
$ws.Range("I66").Value = 22000
$ws.Range("E70").Value = 259600
$ws.Range("F70").Value = 160300
$ws.Range("G70").Value = 71800
$ws.Range("H70").Value = 23400
$ws.Range("I70").Value = 11800
$ws.Range("D72").Value = -212600
$ws.Range("E72").Value = -202300
$ws.Range("F72").Value = -109100
$ws.Range("G72").Value = -44700
$ws.Range("H72").Value = -10400
$ws.Range("I72").Value = -4200
$ws.Range("D76").Value = 190400
$ws.Range("E76").Value = -215200
$ws.Range("F76").Value = -113200
$ws.Range("G76").Value = -44300
$ws.Range("H76").Value = -10000
$ws.Range("I76").Value = -4200
$ws.Range("D81").Value = -10300
$ws.Range("E81").Value = -95000
$ws.Range("F81").Value = -64700
$ws.Range("G81").Value = -34500
$ws.Range("H81").Value = -6400
$ws.Range("I81").Value = -3600
$ws.Range("E83").Value = 2000
$ws.Range("D89").Value = -26300
$ws.Range("E89").Value = -37200
$ws.Range("F89").Value = -18800
$ws.Range("G89").Value = -23900
$ws.Range("H89").Value = -5000
$ws.Range("I89").Value = -4900
$ws.Range("D91").Value = -2900
$ws.Range("F91").Value = -2300
$ws.Range("D94").Value = -46200
$ws.Range("F94").Value = -2300
$ws.Range("G94").Value = -17400
$ws.Range("D100").Value = 133300
$ws.Range("E100").Value = 6600
$ws.Range("F100").Value = 54200
$ws.Range("G100").Value = 44400
$ws.Range("H100").Value = 10800
$ws.Range("I100").Value = 5900
$ws.Range("D101").Value = -1800
$ws.Range("D102").Value = 59000
$ws.Range("E102").Value = -34000
$ws.Range("F102").Value = 31600
$ws.Range("G102").Value = 3200
$ws.Range("H102").Value = 5200
